# Scheduled market-data refresh: update the FFXIV Leve profit sheets
# (currentAveragePrice*/LevePrice*/LeveProfit* columns H:N) with the
# latest pulled prices for the affected Leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 96969
$ws.Range("J3").Value = 96969
$ws.Range("L3").Value = 96969
$ws.Range("N3").Value = -97197

$ws.Range("H34").Value = 4144.75
$ws.Range("I34").Value = 4144.75
$ws.Range("K34").Value = 4144.75
$ws.Range("M34").Value = -3941.75

$ws.Range("H36").Value = 4144.75
$ws.Range("I36").Value = 4144.75
$ws.Range("K36").Value = 4144.75
$ws.Range("M36").Value = -3429.75

$ws.Range("H92").Value = 952.36365
$ws.Range("I92").Value = 867.6
$ws.Range("J92").Value = 1800
$ws.Range("K92").Value = 867.6
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = 380.4
$ws.Range("N92").Value = -4296

$ws.Range("H100").Value = 3595.5557
$ws.Range("I100").Value = 4337.143
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 4337.143
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -3796.143
$ws.Range("N100").Value = -2082

$ws.Range("H102").Value = 96969
$ws.Range("J102").Value = 96969
$ws.Range("L102").Value = 96969
$ws.Range("N102").Value = -103459

$ws.Range("H132").Value = 1798.5385
$ws.Range("I132").Value = 1798.5385
$ws.Range("K132").Value = 5395.6155
$ws.Range("M132").Value = -2865.6155

$ws.Range("H137").Value = 1886.75
$ws.Range("I137").Value = 1823.7826
$ws.Range("K137").Value = 5471.3478
$ws.Range("M137").Value = -2921.3478

$ws.Range("H138").Value = 3501.8333
$ws.Range("J138").Value = 4028.375
$ws.Range("L138").Value = 12085.125
$ws.Range("N138").Value = -22365.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 74500
$ws.Range("J24").Value = 74500
$ws.Range("L24").Value = 74500
$ws.Range("N24").Value = -75248

$ws.Range("H74").Value = 9481.066000000001
$ws.Range("I74").Value = 10225.667
$ws.Range("K74").Value = 10225.667
$ws.Range("M74").Value = -9351.666999999999

$ws.Range("H77").Value = 9481.066000000001
$ws.Range("I77").Value = 10225.667
$ws.Range("K77").Value = 51128.335
$ws.Range("M77").Value = -46760.335

$ws.Range("H100").Value = 74500
$ws.Range("J100").Value = 74500
$ws.Range("L100").Value = 74500
$ws.Range("N100").Value = -76664

$ws.Range("H102").Value = 1173.4667
$ws.Range("I102").Value = 1364.5834
$ws.Range("J102").Value = 409
$ws.Range("K102").Value = 1364.5834
$ws.Range("L102").Value = 409
$ws.Range("M102").Value = 257.4166
$ws.Range("N102").Value = -3653

$ws.Range("H122").Value = 3900
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -15400

$ws.Range("H131").Value = 50650
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws.Range("H132").Value = 3366.9092
$ws.Range("I132").Value = 2735.8572
$ws.Range("J132").Value = 4471.25
$ws.Range("K132").Value = 8207.571599999999
$ws.Range("L132").Value = 13413.75
$ws.Range("M132").Value = -5677.571599999999
$ws.Range("N132").Value = -18473.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4308.5
$ws.Range("I20").Value = 1170.4
$ws.Range("K20").Value = 1170.4
$ws.Range("M20").Value = -923.4000000000001

$ws.Range("H134").Value = 6424.125
$ws.Range("I134").Value = 3565.5
$ws.Range("K134").Value = 10696.5
$ws.Range("M134").Value = -8161.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1809.5
$ws.Range("I31").Value = 1712
$ws.Range("J31").Value = 2199.5
$ws.Range("K31").Value = 1712
$ws.Range("L31").Value = 2199.5
$ws.Range("M31").Value = -1417
$ws.Range("N31").Value = -2789.5

$ws.Range("H34").Value = 1809.5
$ws.Range("I34").Value = 1712
$ws.Range("J34").Value = 2199.5
$ws.Range("K34").Value = 1712
$ws.Range("L34").Value = 2199.5
$ws.Range("M34").Value = -1510
$ws.Range("N34").Value = -2603.5

$ws.Range("H132").Value = 3416.5334
$ws.Range("I132").Value = 2877.9092
$ws.Range("K132").Value = 8633.7276
$ws.Range("M132").Value = -6103.7276

$ws.Range("H134").Value = 10492.167
$ws.Range("I134").Value = 10743.5
$ws.Range("K134").Value = 32230.5
$ws.Range("M134").Value = -29695.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2777.7778
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 9000
$ws.Range("N4").Value = -9224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4081.7778
$ws.Range("I113").Value = 3947
$ws.Range("J113").Value = 4250.25
$ws.Range("K113").Value = 3947
$ws.Range("L113").Value = 4250.25
$ws.Range("M113").Value = -1777
$ws.Range("N113").Value = -8590.25

$ws.Range("H132").Value = 3782.7144
$ws.Range("I132").Value = 3120.5
$ws.Range("K132").Value = 9361.5
$ws.Range("M132").Value = -6831.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 74884
$ws.Range("J76").Value = 74884
$ws.Range("L76").Value = 74884
$ws.Range("N76").Value = -75560

$ws.Range("H79").Value = 74884
$ws.Range("J79").Value = 74884
$ws.Range("L79").Value = 74884
$ws.Range("N79").Value = -77224

$ws.Range("H100").Value = 5316.6665
$ws.Range("I100").Value = 5316.6665
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5316.6665
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = -4775.6665

$ws.Range("H104").Value = 49790
$ws.Range("J104").Value = 49790
$ws.Range("L104").Value = 49790
$ws.Range("N104").Value = -56778

$ws.Range("H106").Value = 55465.668
$ws.Range("J106").Value = 55465.668
$ws.Range("L106").Value = 55465.668
$ws.Range("N106").Value = -57989.668

$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180

$ws.Range("H136").Value = 3115.889
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4739.25
$ws.Range("I81").Value = 1844.8572
$ws.Range("J81").Value = 25000
$ws.Range("K81").Value = 3689.7144
$ws.Range("L81").Value = 50000
$ws.Range("M81").Value = -2628.7144
$ws.Range("N81").Value = -52122

$ws.Range("H84").Value = 4739.25
$ws.Range("I84").Value = 1844.8572
$ws.Range("J84").Value = 25000
$ws.Range("K84").Value = 18448.572
$ws.Range("L84").Value = 250000
$ws.Range("M84").Value = -13144.572
$ws.Range("N84").Value = -260608

$ws.Range("H126").Value = 1808.8
$ws.Range("I126").Value = 1862.25
$ws.Range("K126").Value = 5586.75
$ws.Range("M126").Value = -3116.75

$ws.Range("H130").Value = 89999
$ws.Range("J130").Value = 89999
$ws.Range("L130").Value = 89999
$ws.Range("N130").Value = -100039
